# inicia insercao da tipologia de alteracoes gsim
#
# The existing "Sheet1" (execucao detail rows) gets a new "uo"-level pivot
# summary inserted as a brand-new first sheet, named "Sheet1 (2)" (as Excel
# names a copy made via "Move or Copy > Create a copy" placed before the
# original). The original "Sheet1" also gains two new header columns
# (2009 / 2010) in I1:J1.

$wb = $excel.ActiveWorkbook
$oldWs = $wb.Worksheets.Item("Sheet1")

# New sheet, inserted before the existing "Sheet1" -> becomes sheet1.xml,
# sheetId 2, rId1; old "Sheet1" shifts to sheetId 1 / rId2 / sheet2.xml.
$newWs = $wb.Worksheets.Add($oldWs)
$newWs.Name = "Sheet1 (2)"

# NOTE: worksheet handles obtained via Item(...) track *position*, not the
# sheet identity, so after inserting a new sheet before it, $oldWs now
# actually points at the newly inserted sheet. Re-resolve the original
# "Sheet1" by name so subsequent edits land on the right part.
$oldWs = $wb.Worksheets.Item("Sheet1")

$comma = "_-* #,##0.00_-;-* #,##0.00_-;_-* ""-""??_-;_-@_-"

# ---- formatting first, so values set afterwards inherit the cell style ----

# Header row: centered (reuses the workbook's existing "center" style)
$newWs.Range("B1:I1").HorizontalAlignment = -4108

# Comma-styled columns (reuses the workbook's existing "Comma" style)
$newWs.Range("A2:A9").NumberFormat = $comma
$newWs.Range("F2:H8").NumberFormat = $comma
$newWs.Range("F9:F22").NumberFormat = $comma

$newWs.Columns.Item(5).ColumnWidth = 22

# ---- header values ----
$newWs.Range("B1").Value = "uo_cod"
$newWs.Range("C1").Value = "uo_sigla"
$newWs.Range("D1").Value = "cod"
$newWs.Range("E1").Value = "hist"
$newWs.Range("F1").Value = 2009
$newWs.Range("G1").Value = 2010
$newWs.Range("H1").Value = 2011
$newWs.Range("I1").Value = 2012

# ---- uo / cod / hist pivot rows (2-8) ----
$newWs.Range("B2").Value = 9901
$newWs.Range("C2").Value = "TESOURO"
$newWs.Range("D2").Value = 10
$newWs.Range("E2").Value = "ICMS"
$newWs.Range("F2").Value = 50
$newWs.Range("G2").Value = 50
$newWs.Range("H2").Value = 50

$newWs.Range("B3").Value = 9901
$newWs.Range("C3").Value = "TESOURO"
$newWs.Range("D3").Value = 20
$newWs.Range("E3").Value = "Taxa de segurança pública"
$newWs.Range("F3").Value = 25
$newWs.Range("G3").Value = 25
$newWs.Range("H3").Value = 25

$newWs.Range("B4").Value = 2301
$newWs.Range("C4").Value = "DER"
$newWs.Range("D4").Value = 20
$newWs.Range("E4").Value = "CIDE"
$newWs.Range("F4").Value = 15
$newWs.Range("G4").Value = 15
$newWs.Range("H4").Value = 15

$newWs.Range("B5").Value = 1251
$newWs.Range("C5").Value = "PMMG"
$newWs.Range("D5").Value = 60
$newWs.Range("E5").Value = "Concursos"
$newWs.Range("F5").Value = 10
$newWs.Range("G5").Value = 0
$newWs.Range("H5").Value = 10

$newWs.Range("B6").Value = 1251
$newWs.Range("C6").Value = "PMMG"
$newWs.Range("D6").Value = 10
$newWs.Range("E6").Value = "Concursos"
$newWs.Range("F6").Value = 0
$newWs.Range("G6").Value = 10
$newWs.Range("H6").Value = 0

$newWs.Range("B7").Value = 2271
$newWs.Range("C7").Value = "FHEMIG"
$newWs.Range("D7").Value = 60
$newWs.Range("E7").Value = "SUS"
$newWs.Range("F7").Value = 5
$newWs.Range("G7").Value = 0
$newWs.Range("H7").Value = 0

$newWs.Range("B8").Value = 2271
$newWs.Range("C8").Value = "FHEMIG"
$newWs.Range("D8").Value = 20
$newWs.Range("E8").Value = "SUS"
$newWs.Range("F8").Value = 0
$newWs.Range("G8").Value = 5
$newWs.Range("H8").Value = 5

# ---- leftover scratch block (rows 18-22), mirrors the first 5 detail
#      rows of the original "Sheet1" pasted below the pivot table ----
$newWs.Range("C18").Value = "TESOURO"
$newWs.Range("D18").Value = 10
$newWs.Range("E18").Value = "ICMS"
$newWs.Range("F18").Value = 50

$newWs.Range("C19").Value = "TESOURO"
$newWs.Range("D19").Value = 20
$newWs.Range("E19").Value = "Taxa de segurança pública"
$newWs.Range("F19").Value = 25

$newWs.Range("C20").Value = "DER"
$newWs.Range("D20").Value = 20
$newWs.Range("E20").Value = "CIDE"
$newWs.Range("F20").Value = 15

$newWs.Range("C21").Value = "PMMG"
$newWs.Range("D21").Value = 60
$newWs.Range("E21").Value = "Concursos"
$newWs.Range("F21").Value = 10
$newWs.Range("G21").Value = "x"

$newWs.Range("C22").Value = "FHEMIG"
$newWs.Range("D22").Value = 20
$newWs.Range("E22").Value = "SUS"
$newWs.Range("F22").Value = 5
$newWs.Range("G22").Value = "FHEMIG não voltou para a fonte 60"

# Active selection on the new sheet
$newWs.Range("H8").Select()

# ---- original "Sheet1": two new trailing year columns on the header row ----
$oldWs.Range("I1").Value = 2009
$oldWs.Range("J1").Value = 2010
